$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same header formatting used by the existing "07-04-2025" columns
# (D1:E1) to the two new date-column pairs being added (F1:G1 and H1:I1).
$ws.Range("D1:E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1:I1").PasteSpecial(-4122) # xlPasteFormats

# New header labels for 09-04-2025 and 10-04-2025
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Default attendance data for every student row: absent status + zero time
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 6).Value = "A"
    $ws.Cells.Item($row, 7).Value = "00:00:00"
    $ws.Cells.Item($row, 8).Value = "A"
    $ws.Cells.Item($row, 9).Value = "00:00:00"
}
